# Atualização dos docentes credenciados.
# Insert two new faculty "header" rows (João Felipe Nicolaci Pimentel and
# Pedro Cortez Fetter Lopes) into the alphabetically-sorted list on Sheet1,
# each styled like the pre-existing "newly added faculty" header rows
# (bigger/grey font on the name, text-number-format on the Lattes/Scholar
# ids) but additionally with a taller row (16pt) and a distinct font.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 24: João Felipe Nicolaci Pimentel (inserted before José Viterbo Filho) ---
$ws.Range("A24").EntireRow.Insert()
$ws.Rows("24:24").RowHeight() = 16

$ws.Range("B24").NumberFormat() = "@"
$ws.Range("C24").NumberFormat() = "@"

$ws.Range("A24").Value() = "João Felipe Nicolaci Pimentel"
$ws.Range("B24").Value() = "4603791761884563"
$ws.Range("C24").Value() = "1rMhTTcAAAAJ"

# --- Row 43: Pedro Cortez Fetter Lopes (inserted before Raphael Carlos Santos Machado) ---
$ws.Range("A43").EntireRow.Insert()
$ws.Rows("43:43").RowHeight() = 16

$ws.Range("B43").NumberFormat() = "@"
$ws.Range("C43").NumberFormat() = "@"

$ws.Range("A43").Value() = "Pedro Cortez Fetter Lopes"
$ws.Range("B43").Value() = "7336253957211512"
$ws.Range("C43").Value() = "qDSsZdgAAAAJ"

# Font for both new "name" cells: 12pt, grey (#555555). The id cells
# (B/C) keep the plain text-number-format style already used by the other
# "new faculty" rows.
$ws.Range("A24").Font.Size() = 12
$ws.Range("A24").Font.Color() = 5592405
$ws.Range("A43").Font.Size() = 12
$ws.Range("A43").Font.Color() = 5592405

# Selection moved off the inserted rows, matching the saved file.
$ws.Range("C44").Select()
